# The workbook's single sheet is being renamed from "Síntese" to "Sintese"
# (accent removed) as part of the export-folder restructuring described in
# the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Sintese"
